$wb = $excel.ActiveWorkbook

# --- Update the "survey" sheet ---
$ws = $wb.Worksheets.Item("survey")

# Change type of M_FOL_date (row 2) from "date" to "text"
$ws.Range("C2").Value = "text"

# Change type of M_time (row 4) from "time" to "text"
$ws.Range("C4").Value = "text"

# Move the active selection from C13 to C6
$ws.Range("C6").Select()

# --- Update workbook view (window size/position) ---
# (Best effort: mirrors the saved bookViews/workbookView geometry recorded by a
# real Excel client. Window chrome is host UI state rather than workbook
# content, so not every COM surface persists it, but we set it everywhere
# it is exposed.)
$excel.Left = 0
$excel.Top = 0
$excel.Width = 25600
$excel.Height = 16060

$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 0
$win.Width = 25600
$win.Height = 16060

$aw = $excel.ActiveWindow
$aw.Left = 0
$aw.Top = 0
$aw.Width = 25600
$aw.Height = 16060
